# Add spinner UI chrome to show when results are loading.
#
# Semantic changes on Sheet1 ("tpivot next steps" tracker):
#   - Row 2 ("Show UI chrome while parse request is being processed")
#     status moves from "Not started" to "Complete".
#   - Row 9 ("Refactor Jquery code...") gets its Approach note reworded
#     (now that the view/model split has landed) and its status flips to
#     "Complete" too; since it no longer matches the "Not started"
#     AutoFilter criterion it is hidden.
#   - Selection moves to F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: Status -> Complete
$ws.Range("F2").Value = "Complete"

# Row 9: reword Approach note, flip status to Complete, and hide the row
$ws.Range("E9").Value = "View/model operations are already lexically separate. Not difficult to move into own modules."
$ws.Range("F9").Value = "Complete"
$ws.Rows.Item(9).Hidden = $true

# Move the active selection to F7
$ws.Range("F7").Select()
